$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()
$wb.Application.ActiveWindow.ScrollRow = 7
$wb.Application.ActiveWindow.ScrollColumn = 1
$ws3.Range("B13").Select()
Write-Output "done"
